$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Update cell H6 first (so its new shared string is registered before G6's),
# replacing the Screenshot validations with isIconDisplayed validations
$h6 = "validate1`n{`nvalidate_PageTitle=Compliance JS specs`n};`nvalidate2`n{`nvalidate_PageTitle=Native Toolbar JS Test`n};`nvalidate3`n{`nvalidate_Text_Exists=VT200-0603`n};`nvalidate4`n{`nvalidate_isIconDisplayed=toobarview_xpath,true`n};`nvalidate5`n{`nvalidate_Text_Exists=VT200-0605`n};`nvalidate6`n{`nvalidate_isIconDisplayed=toobarview_xpath,false`n};`n"
$ws.Range("H6").Value2 = $h6

# Update cell G6: remove the two TakeScreenshot(...) lines from the step script
$g6 = "wait(5);`nvalidate1;`nlink_Click(toolbar_test_link);`nvalidate2;`nSelectTestToRun(VT200_0603_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nClickRunTest(runtest_bottom_xpath);`nwait(5);`nvalidate4;`nSelectTestToRun(VT200_0605_string);`nClickRunTest(runtest_top_xpath);`nvalidate5;`nClickRunTest(runtest_bottom_xpath);`nwait(3);`nvalidate6;"
$ws.Range("G6").Value2 = $g6

# Update sheet view/selection: clear the scrolled topLeftCell and select E1 instead of G6
$ws.Activate()
$ws.Range("E1").Select()
